# Update Leve profit/price figures per scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 498.77777
$ws.Range("I32").Value = 425
$ws.Range("J32").Value = 519.8570999999999
$ws.Range("K32").Value = 425
$ws.Range("L32").Value = 519.8570999999999
$ws.Range("M32").Value = -99
$ws.Range("N32").Value = -1171.8571

$ws.Range("H99").Value = 457.2
$ws.Range("I99").Value = 228
$ws.Range("J99").Value = 882.8570999999999
$ws.Range("K99").Value = 684
$ws.Range("L99").Value = 2648.5713
$ws.Range("M99").Value = 814
$ws.Range("N99").Value = -5644.5713

$ws.Range("H137").Value = 1726.2593
$ws.Range("I137").Value = 1548.0476
$ws.Range("K137").Value = 4644.142800000001
$ws.Range("M137").Value = -2094.142800000001

$ws.Range("H138").Value = 1945.258
$ws.Range("I138").Value = 1690.5
$ws.Range("J138").Value = 2106.158
$ws.Range("K138").Value = 5071.5
$ws.Range("L138").Value = 6318.474
$ws.Range("M138").Value = 68.5
$ws.Range("N138").Value = -16598.474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 6050
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 6050
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 6050
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = -7414

$ws.Range("H61").Value = 1798.875
$ws.Range("I61").Value = 798.36365
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 798.36365
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -586.36365
$ws.Range("N61").Value = -4424

$ws.Range("H88").Value = 459253.03
$ws.Range("I88").Value = 717990
$ws.Range("J88").Value = 6463.375
$ws.Range("K88").Value = 717990
$ws.Range("L88").Value = 6463.375
$ws.Range("M88").Value = -717584
$ws.Range("N88").Value = -7275.375

$ws.Range("H91").Value = 459253.03
$ws.Range("I91").Value = 717990
$ws.Range("J91").Value = 6463.375
$ws.Range("K91").Value = 717990
$ws.Range("L91").Value = 6463.375
$ws.Range("M91").Value = -716586
$ws.Range("N91").Value = -9271.375

$ws.Range("H136").Value = 1798.875
$ws.Range("I136").Value = 798.36365
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2395.09095
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 154.9090500000002
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 9666.666999999999
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4828

$ws.Range("H86").Value = 1749.4
$ws.Range("I86").Value = 1841
$ws.Range("J86").Value = 1535.6666
$ws.Range("K86").Value = 1841
$ws.Range("L86").Value = 1535.6666
$ws.Range("M86").Value = -718
$ws.Range("N86").Value = -3781.6666

$ws.Range("H89").Value = 1749.4
$ws.Range("I89").Value = 1841
$ws.Range("J89").Value = 1535.6666
$ws.Range("K89").Value = 9205
$ws.Range("L89").Value = 7678.333000000001
$ws.Range("M89").Value = -3589
$ws.Range("N89").Value = -18910.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6670233
$ws.Range("I31").Value = 2481.2222
$ws.Range("K31").Value = 2481.2222
$ws.Range("M31").Value = -2186.2222

$ws.Range("H34").Value = 6670233
$ws.Range("I34").Value = 2481.2222
$ws.Range("K34").Value = 2481.2222
$ws.Range("M34").Value = -2279.2222

$ws.Range("H62").Value = 66670268
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 83336830
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 83336830
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -83338078

$ws.Range("H65").Value = 66670268
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 83336830
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 416684150
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -416690390

$ws.Range("H132").Value = 2957.1765
$ws.Range("I132").Value = 2022.8334
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 6068.5002
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -3538.5002
$ws.Range("N132").Value = -20658.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2828.5715
$ws.Range("I80").Value = 2650
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 7950
$ws.Range("L80").Value = 8700
$ws.Range("M80").Value = -7014
$ws.Range("N80").Value = -10572

$ws.Range("H83").Value = 2828.5715
$ws.Range("I83").Value = 2650
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 23850
$ws.Range("L83").Value = 26100
$ws.Range("M83").Value = -19170
$ws.Range("N83").Value = -35460

$ws.Range("H131").Value = 3972551
$ws.Range("I131").Value = 8745
$ws.Range("J131").Value = 6945405.5
$ws.Range("K131").Value = 26235
$ws.Range("L131").Value = 20836216.5
$ws.Range("M131").Value = -21195
$ws.Range("N131").Value = -20846296.5

$ws.Range("H139").Value = 3165.7144
$ws.Range("I139").Value = 1380
$ws.Range("J139").Value = 3880
$ws.Range("K139").Value = 4140
$ws.Range("L139").Value = 11640
$ws.Range("M139").Value = 1000
$ws.Range("N139").Value = -21920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -11288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 6515
$ws.Range("I32").Value = 6515
$ws.Range("K32").Value = 6515
$ws.Range("M32").Value = -6198

$ws.Range("H124").Value = 134714.5
$ws.Range("J124").Value = 134714.5
$ws.Range("L124").Value = 134714.5
$ws.Range("N124").Value = -144534.5

$ws.Range("H132").Value = 2139.5715
$ws.Range("I132").Value = 1495.1
$ws.Range("J132").Value = 3750.75
$ws.Range("K132").Value = 4485.299999999999
$ws.Range("L132").Value = 11252.25
$ws.Range("M132").Value = -1955.299999999999
$ws.Range("N132").Value = -16312.25

$ws.Range("H136").Value = 6877.4287
$ws.Range("I136").Value = 8042.706
$ws.Range("J136").Value = 1925
$ws.Range("K136").Value = 24128.118
$ws.Range("L136").Value = 5775
$ws.Range("M136").Value = -21578.118
$ws.Range("N136").Value = -10875
